$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B5" = 13.4133266137753
    "C5" = 0.0016615388127017398
    "D5" = 0.03551512420339932
    "B7" = 103.04134659209438
    "C7" = 0.0089688427267222339
    "D7" = 5.8512140324020248
    "B8" = 107.18904048327813
    "C8" = 0.0057050094697372001
    "D8" = 1.2620154565463195
    "B11" = 15.172265319927243
    "C11" = 0.00057755050555109591
    "D11" = 0.70038464050956661
    "B12" = 169.20631933846897
    "C12" = 0.0098356698275528977
    "D12" = 3.4618629466157214
    "B14" = 128.32336925006732
    "C14" = 0.0072670020744946429
    "D14" = 1.0440624022049168
    "B17" = 22.200185634078924
    "C17" = 0.018423154879788857
    "D17" = 0.47308634515264303
    "B18" = 161.09476130796733
    "C18" = 0.0041840544910831182
    "D18" = 7.0700967237073886
    "B19" = 156.28493648011718
    "C19" = 0.0051729385348202482
    "D19" = 5.9668455055760061
    "B20" = 1179.5102111569222
    "C20" = 0.099132876628876571
    "D20" = 18.936228907481972
    "B21" = 16.79694362006747
    "C21" = 0.00096293317599829072
    "D21" = 0.29446627408051268
    "B22" = 9.8566313313667386
    "C22" = 0.00030992328239202337
    "D22" = 0.82557572614867492
    "B23" = 157.93570949823342
    "C23" = 0.046158647254689031
    "D23" = 4.5681723753020043
    "B24" = 69.652735435237545
    "C24" = 0.004535140423012448
    "D24" = 0.077434068268624243
    "B25" = 81.8264172319073
    "C25" = 0.0012533746715075959
    "D25" = 6.417859515720628
    "B26" = 416.44570418955772
    "C26" = 0.054397597041174431
    "D26" = 32.459020443473008
    "B27" = 375.44714188058316
    "C27" = 0.028809118384815063
    "D27" = 3.039775168031901
    "B28" = 216.2911153880159
    "C28" = 0.012366037097810784
    "D28" = 0.52955117152995668
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
